# Regenerate orders with updated distance/size codes.
#
# The underlying shared strings encode trial conditions such as
# "Face11_D64_S25" / "Face11_D64_S25_l.png" / "D64" / "S30", where the
# "D.." token is a viewing distance and the "S.." token is a stimulus size.
# This edit renumbers those codes (affecting every column built from them:
# Condition, Filename_Left, Filename_Right, Distance, Size) via a simple
# substring replacement, applied uniformly to every text cell in the sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = $ws.UsedRange.Rows.Count()
$cols = $ws.UsedRange.Columns.Count()

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            $nv = $v.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
